# Excel COM-interop script.
#
# Commit: "Generovani ok predmetu a zaklad pro xlsx -> csv"
#
# Net effect observed in the target OOXML diff: the two data rows for
# r14 (katedra=KBI, predmet="Evolucni biologie") and r15
# (katedra=KFY, predmet="PVK-Kvantova fyzika II") trade places - every
# column's content on row 14 becomes what row 15 used to hold and vice
# versa (the shared-string-table renumbering in the raw XML is just a
# side effect of that row swap, not an independent change).
#
# We therefore write out the literal (pre-edit) content of each row into
# the other row. Values are hard-coded (rather than read back with
# `.Value`, which this host's COM shim cannot round-trip through a
# variable) and written through here-strings so none of the embedded
# quotes/newlines in the long free-text columns need escaping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 14 <- original row 15 content ---
$n14_A = @'
KFY
'@
$ws.Cells.Item(14, 1).Value = $n14_A

$n14_B = @'
K607
'@
$ws.Cells.Item(14, 2).Value = $n14_B

$n14_C = 2023
$ws.Cells.Item(14, 3).Value = $n14_C

$n14_D = @'
PVK-Kvantová fyzika II
'@
$ws.Cells.Item(14, 4).Value = $n14_D

$n14_E = @'
PVK-Kvantová fyzika II
'@
$ws.Cells.Item(14, 5).Value = $n14_E

$n14_F = @'
N
'@
$ws.Cells.Item(14, 6).Value = $n14_F

$n14_G = @'
A
'@
$ws.Cells.Item(14, 7).Value = $n14_G

$n14_H = 5
$ws.Cells.Item(14, 8).Value = $n14_H

$n14_I = @'
NE
'@
$ws.Cells.Item(14, 9).Value = $n14_I

$n14_J = $null
$ws.Cells.Item(14, 10).Value = $n14_J

$n14_K = $null
$ws.Cells.Item(14, 11).Value = $n14_K

$n14_L = $null
$ws.Cells.Item(14, 12).Value = $n14_L

$n14_M = @'
None
'@
$ws.Cells.Item(14, 13).Value = $n14_M

$n14_N = $null
$ws.Cells.Item(14, 14).Value = $n14_N

$n14_O = $null
$ws.Cells.Item(14, 15).Value = $n14_O

$n14_P = @'
None
'@
$ws.Cells.Item(14, 16).Value = $n14_P

$n14_Q = $null
$ws.Cells.Item(14, 17).Value = $n14_Q

$n14_R = $null
$ws.Cells.Item(14, 18).Value = $n14_R

$n14_S = @'
None
'@
$ws.Cells.Item(14, 19).Value = $n14_S

$n14_T = $null
$ws.Cells.Item(14, 20).Value = $n14_T

$n14_U = $null
$ws.Cells.Item(14, 21).Value = $n14_U

$n14_V = @'
None
'@
$ws.Cells.Item(14, 22).Value = $n14_V

$n14_W = $null
$ws.Cells.Item(14, 23).Value = $n14_W

$n14_X = $null
$ws.Cells.Item(14, 24).Value = $n14_X

$n14_Y = $null
$ws.Cells.Item(14, 25).Value = $n14_Y

$n14_Z = $null
$ws.Cells.Item(14, 26).Value = $n14_Z

$n14_AA = $null
$ws.Cells.Item(14, 27).Value = $n14_AA

$n14_AB = $null
$ws.Cells.Item(14, 28).Value = $n14_AB

$n14_AC = $null
$ws.Cells.Item(14, 29).Value = $n14_AC

$n14_AD = $null
$ws.Cells.Item(14, 30).Value = $n14_AD

$n14_AE = $null
$ws.Cells.Item(14, 31).Value = $n14_AE

$n14_AF = $null
$ws.Cells.Item(14, 32).Value = $n14_AF

$n14_AG = $null
$ws.Cells.Item(14, 33).Value = $n14_AG

$n14_AH = @'
A
'@
$ws.Cells.Item(14, 34).Value = $n14_AH

$n14_AI = 12
$ws.Cells.Item(14, 35).Value = $n14_AI

$n14_AJ = @'
HOD/TYD
'@
$ws.Cells.Item(14, 36).Value = $n14_AJ

$n14_AK = 4
$ws.Cells.Item(14, 37).Value = $n14_AK

$n14_AL = @'
HOD/TYD
'@
$ws.Cells.Item(14, 38).Value = $n14_AL

$n14_AM = 0
$ws.Cells.Item(14, 39).Value = $n14_AM

$n14_AN = @'
HOD/TYD
'@
$ws.Cells.Item(14, 40).Value = $n14_AN

$n14_AO = $null
$ws.Cells.Item(14, 41).Value = $n14_AO

$n14_AP = @'
Zkouška
'@
$ws.Cells.Item(14, 42).Value = $n14_AP

$n14_AQ = @'
ANO
'@
$ws.Cells.Item(14, 43).Value = $n14_AQ

$n14_AR = @'
Kombinovaná
'@
$ws.Cells.Item(14, 44).Value = $n14_AR

$n14_AS = $null
$ws.Cells.Item(14, 45).Value = $n14_AS

$n14_AT = $null
$ws.Cells.Item(14, 46).Value = $n14_AT

$n14_AU = $null
$ws.Cells.Item(14, 47).Value = $n14_AU

$n14_AV = $null
$ws.Cells.Item(14, 48).Value = $n14_AV

$n14_AW = $null
$ws.Cells.Item(14, 49).Value = $n14_AW

$n14_AX = $null
$ws.Cells.Item(14, 50).Value = $n14_AX

$n14_AY = @'
Čeština
'@
$ws.Cells.Item(14, 51).Value = $n14_AY

$n14_AZ = $null
$ws.Cells.Item(14, 52).Value = $n14_AZ

$n14_BA = @'
A
'@
$ws.Cells.Item(14, 53).Value = $n14_BA

$n14_BB = @'
N
'@
$ws.Cells.Item(14, 54).Value = $n14_BB

$n14_BC = @'
N
'@
$ws.Cells.Item(14, 55).Value = $n14_BC

$n14_BD = $null
$ws.Cells.Item(14, 56).Value = $n14_BD

$n14_BE = $null
$ws.Cells.Item(14, 57).Value = $n14_BE

$n14_BF = $null
$ws.Cells.Item(14, 58).Value = $n14_BF

$n14_BG = @'
N
'@
$ws.Cells.Item(14, 59).Value = $n14_BG

$n14_BH = @'
N
'@
$ws.Cells.Item(14, 60).Value = $n14_BH

$n14_BI = $null
$ws.Cells.Item(14, 61).Value = $n14_BI

$n14_BJ = 0
$ws.Cells.Item(14, 62).Value = $n14_BJ

$n14_BK = @'
Bc.
'@
$ws.Cells.Item(14, 63).Value = $n14_BK

$n14_BL = @'
Mgr.
'@
$ws.Cells.Item(14, 64).Value = $n14_BL

$n14_BM = @'
A
'@
$ws.Cells.Item(14, 65).Value = $n14_BM

$n14_BN = @'
None
'@
$ws.Cells.Item(14, 66).Value = $n14_BN

$n14_BO = $null
$ws.Cells.Item(14, 67).Value = $n14_BO

$n14_BP = $null
$ws.Cells.Item(14, 68).Value = $n14_BP

$n14_BQ = $null
$ws.Cells.Item(14, 69).Value = $n14_BQ

# --- row 15 <- original row 14 content ---
$n15_A = @'
KBI
'@
$ws.Cells.Item(15, 1).Value = $n15_A

$n15_B = @'
K607
'@
$ws.Cells.Item(15, 2).Value = $n15_B

$n15_C = 2023
$ws.Cells.Item(15, 3).Value = $n15_C

$n15_D = @'
Evoluční biologie
'@
$ws.Cells.Item(15, 4).Value = $n15_D

$n15_E = @'
Evoluční biologie
'@
$ws.Cells.Item(15, 5).Value = $n15_E

$n15_F = @'
N
'@
$ws.Cells.Item(15, 6).Value = $n15_F

$n15_G = @'
A
'@
$ws.Cells.Item(15, 7).Value = $n15_G

$n15_H = 5
$ws.Cells.Item(15, 8).Value = $n15_H

$n15_I = @'
NE
'@
$ws.Cells.Item(15, 9).Value = $n15_I

$n15_J = $null
$ws.Cells.Item(15, 10).Value = $n15_J

$n15_K = @'
''RNDr. Jan Ipser, CSc.'
'@
$ws.Cells.Item(15, 11).Value = $n15_K

$n15_L = @'
''RNDr. Jan Ipser, CSc.' (100)
'@
$ws.Cells.Item(15, 12).Value = $n15_L

$n15_M = @'
[930]
'@
$ws.Cells.Item(15, 13).Value = $n15_M

$n15_N = $null
$ws.Cells.Item(15, 14).Value = $n15_N

$n15_O = $null
$ws.Cells.Item(15, 15).Value = $n15_O

$n15_P = @'
None
'@
$ws.Cells.Item(15, 16).Value = $n15_P

$n15_Q = $null
$ws.Cells.Item(15, 17).Value = $n15_Q

$n15_R = $null
$ws.Cells.Item(15, 18).Value = $n15_R

$n15_S = @'
None
'@
$ws.Cells.Item(15, 19).Value = $n15_S

$n15_T = $null
$ws.Cells.Item(15, 20).Value = $n15_T

$n15_U = $null
$ws.Cells.Item(15, 21).Value = $n15_U

$n15_V = @'
None
'@
$ws.Cells.Item(15, 22).Value = $n15_V

$n15_W = $null
$ws.Cells.Item(15, 23).Value = $n15_W

$n15_X = $null
$ws.Cells.Item(15, 24).Value = $n15_X

$n15_Y = $null
$ws.Cells.Item(15, 25).Value = $n15_Y

$n15_Z = $null
$ws.Cells.Item(15, 26).Value = $n15_Z

$n15_AA = @'
KBI/K101
'@
$ws.Cells.Item(15, 27).Value = $n15_AA

$n15_AB = $null
$ws.Cells.Item(15, 28).Value = $n15_AB

$n15_AC = $null
$ws.Cells.Item(15, 29).Value = $n15_AC

$n15_AD = @'
''Flegr J. Evoluční biologie. Academia Praha, 2005. ',
'Flegr J. Evoluční biologie. Academia Praha, 2005. ',
'Rosypal S. a kol. Fylogeneze, systém a biologie organismů. SPN Praha, 1992. ',
'Rosypal S. a kol. Fylogeneze, systém a biologie organismů. SPN Praha, 1992. ',
'Rosypal a kol. Nový přehled biologie. 2003. ',
'Rosypal S. a kol. Nový přehled biologie. 2003. ',
'http://biology.ujep.cz/vyuka',
'Dawkins, R. Boží blud. Academia, Praha 2009. ',
'Mayr, E. Co je evoluce? Academia, Praha 2009. ',
'Ohno, S. Evoluce genovou duplikací. Academia, Praha 1975. ',
'Lovelock, J.:. Gaia: Živoucí planeta, Praha: MF, 1994. ',
'Dostál P. a kol. Historický vývoj organismů. UK Praha, 2004. ',
'Dostál P. a kol. Historický vývoj organismů. UK Praha, 2004. ',
'Dostál P. a kol. Historický vývoj organismů. UK Praha, 2004. ',
'null',
'Zrzavý J., Storch D., Mihulka S. Jak se dělá evoluce. Paseka Praha, 2004. ',
'Zrzavý, J., Storch, D., Mihulka, S. Jak se dělá evoluce. Paseka, Praha 2004. ',
'Vácha, M. O. Návrat ke stromu života - evoluce a křesťanství. Cesta, Brno 2005. ',
'Kovář L. Nevyřešné otazníky evoluce. Rubico, Praha 2003. ',
'Darwin, Ch. O pohlavním výběru. Academia, Praha 2005. ',
'Darwin Ch. O vzniku druhů přírodním výběrem. Praha, 1953. ',
'Markoš A. Povstávání živého tvaru. Vesmír Praha, 1998. ',
'Dawkins, R. Příběh předka. Academia, Praha 2008. ',
'Dunbar R. Příběh rodu Homo. Academia, Praha 2009. ',
'Dawkins R. Sobecký gen. Mladá Fronta Praha. 1998. ',
'Margulisová L. Symbiotická planeta. Academia Praha, 2004. ',
'Margulisová, L. Symbiotická planeta. Academia, Praha 2004. ',
'Moreland, J. P., Reynolds, J. M. Třikrát evoluce versus stvoření. Návrat domů, Praha 2004. ',
'Přívratský, V. Tvar a prostředí v lidské evoluci. Univerzita Karlova, Praha 2003. ',
'Tresmontant, C. Základy teologie. Barrister  Principal, Brno 1995. ',
'Flegr J. Zamrzlá evoluce. Academia Praha, 2006. '
'@
$ws.Cells.Item(15, 30).Value = $n15_AD

$n15_AE = $null
$ws.Cells.Item(15, 31).Value = $n15_AE

$n15_AF = $null
$ws.Cells.Item(15, 32).Value = $n15_AF

$n15_AG = $null
$ws.Cells.Item(15, 33).Value = $n15_AG

$n15_AH = @'
N
'@
$ws.Cells.Item(15, 34).Value = $n15_AH

$n15_AI = 4
$ws.Cells.Item(15, 35).Value = $n15_AI

$n15_AJ = @'
HOD/SEM
'@
$ws.Cells.Item(15, 36).Value = $n15_AJ

$n15_AK = 5
$ws.Cells.Item(15, 37).Value = $n15_AK

$n15_AL = @'
HOD/SEM
'@
$ws.Cells.Item(15, 38).Value = $n15_AL

$n15_AM = 0
$ws.Cells.Item(15, 39).Value = $n15_AM

$n15_AN = @'
HOD/SEM
'@
$ws.Cells.Item(15, 40).Value = $n15_AN

$n15_AO = @'
Cílem předmětu je zprostředkovat studentům informace týkající se biologické evoluce v kontextu s chemickou evolucí a vývojem vesmíru, vztahu evolucionizmu a krecionizmu, historického vývoje evolučních teorií, mikroevoluce a makroevoluce (mechanizmy, principy, metody studia, metodologie a filozofické zobecnění), komparace evolučních a fylogenetických aspektů výsledků paleontologických a molekulárně-biologických výzkumů, aplikace evolučně-biologických poznatků ve společenské praxi.

Tento kurz byl inovován v rámci projektu CZ.1.07/2.2.00/28.0296 "Mezioborové vazby a podpora praxe v přírodovědných a technických studijních programech UJEP"
'@
$ws.Cells.Item(15, 41).Value = $n15_AO

$n15_AP = @'
Zkouška
'@
$ws.Cells.Item(15, 42).Value = $n15_AP

$n15_AQ = @'
NE
'@
$ws.Cells.Item(15, 43).Value = $n15_AQ

$n15_AR = @'
Ústní
'@
$ws.Cells.Item(15, 44).Value = $n15_AR

$n15_AS = @'
Na konzultacích jsou rámcově odpřednášena stěžejní témata (hlavní teorie biologické evoluce, mechanismy evoluce biologických systémů, mikroevoluce a makroevoluce, teorie vzniku života). Následně jsou diskutována a procvičena na cvičení řešením typových úloh k jednotlivým tematickým blokům. Získané informace si posluchači rozšíří samostudiem studijní opory, doporučené odborné literatury a dalších zdrojů. Samostudium představuje hlavní formu přípravy k zápočtu a ke zkoušce.
'@
$ws.Cells.Item(15, 45).Value = $n15_AS

$n15_AT = @'
1. Evoluce jako obecná vlastnost hmoty. Evoluce jako téma vědy, filozofie a teologie. Evoluce     na fyzikální, chemické a biologické úrovni. Vztah evoluce - fylogeneze - protobiologie. 
2. -3. Vývoj evolučního myšlení od antiky po současnost. Charakteristika významných hypotéz a teorií biologické evoluce - Darwinova teorie přírodního výběru, syntetická teorie evoluce, evoluce genovou duplikací (Ohno), teorie sobeckého genu (Dawkins), symbiotická teorie (Margulisová), neutrální teorie (Kimura), téměř neutrální teorie (Ohtová), teorie přerušovaných rovnováh (Gould, Eldridge), zamrzlá evoluce (Flegr). Evolucionizmus versus kreacionizmus; specifické rysy soudobého kreacionizmu.
4. Biologická evoluce - mikroevoluce a makroevoluce. Hierarchie evolučního procesu. Biologická (genetická) variabilita a polymorfismus na různých úrovních biologických systémů (molekulární, buněčná, organizmální, populační, společenstev, druhová). Koevoluce.
5. - 7. Mechanizmy biologické evoluce. 
A. Přírodní a pohlavní výběr v Darwinově teorii, neodarwinizmu (v syntetické teorii evoluce) a v postneodarwinizmu. 
B. Genetický drift.
C. Evoluční tahy (mutační, molekulární, meiotický, reparační)
D. Další evoluční jevy a procesy: exaptace (preadaptace), konvergence (homoplázie), heterochronie, analogie versus homologie, extinkce a pseudoextinkce, evoluční omezení.
8. Místo deterministických procesů (mutace, selekce a migrace) v evoluci biologických systémů. Migrace jednosměrná a obousměrná, migrační rychlost, migrace ve vztahu k izolaci (sub)populace; mutace jako primární událost v evoluci biologických systémů, klasifikace mutací, mutační rychlost; typy selekce, selekční koeficient, adaptivní hodnota a reprodukční zdatnost, interakce mutace a selekce; vliv migrace, mutace a selekce na genetickou strukturu populací; důsledky deterministických (adaptačních) procesů v mikroevoluci.
9. Místo stochastických procesů v evoluci biologických systémů. Disperzivní proces v populacích; genetický drift, inbrídink; efekt zakladatele, efekt hrdla láhve; důsledky stochastických procesů v mikroevoluci.
10 - 11. Speciace. Koncepce druhu (morfologická, fylogenetická, evoluční). Fyletická speciace (změna) versus štěpná (evoluční) speciace; typy a mechanizmy speciace (alopatrická, peripatrická, parapatrická, sympatrická; hybridní zóna, adaptivní radiace; izolace geografická a biologická; divergence a konvergence). Domestikace. 
12. Evoluce na molekulární a buněčné úrovni. Evoluce nukleových kyselin (struktura, funkce) a genetického kódu; RNA svět, DNA svět, ribozymy; vznik a evoluce genetického kódu; evoluce genu; úloha mutace a rekombinace v evoluci - evoluce genovou duplikací (Ohno); vznik a evoluce virů, prokaryotické a eukaryotické buňky; kompartmentace eukaryotické buňky a evoluce celulárních struktur, endosymbiotická teorie (mitochondrie, plastidy); evoluce hlavních metabolických typů a dějů (fotosyntéza, respirace). 
13. Makroevoluce biologických systémů v kontextu s geologickým vývojem Země. Nástin hlavních linií vývoje života od prekambria po kvartér (prekambrium, kambrická exploze, Ediakara, Burgessovy břidlice, interpretace paleontologického materiálu). Fylogenetické stromy - typy, interpretace, metody a postupy jejich konstrukce, využití bioinformatiky. 
14. Protobiologie - vznik života. Hypotézy a teorie vzniku (resp. původu) života na Zemi. Teorie evoluční abiogeneze (naivn abiogeneze - panspermie, Oparinova koacervátová teorie), Millerův a Ureyův pokus, prebiotická polévka, teorie koacervátu v koacervátu (Liebl), teorie mikrosfér (Fox), teorie jílových částic (Bernal, Cairns-Smith), teorie hydrotermálních průduchů; LUCA. 


'@
$ws.Cells.Item(15, 46).Value = $n15_AT

$n15_AU = @'
Absolvované základní kurzy Genetika, Obecná botanika a Obecná zoologie.
'@
$ws.Cells.Item(15, 47).Value = $n15_AU

$n15_AV = $null
$ws.Cells.Item(15, 48).Value = $n15_AV

$n15_AW = $null
$ws.Cells.Item(15, 49).Value = $n15_AW

$n15_AX = $null
$ws.Cells.Item(15, 50).Value = $n15_AX

$n15_AY = @'
Čeština
'@
$ws.Cells.Item(15, 51).Value = $n15_AY

$n15_AZ = $null
$ws.Cells.Item(15, 52).Value = $n15_AZ

$n15_BA = @'
A
'@
$ws.Cells.Item(15, 53).Value = $n15_BA

$n15_BB = @'
N
'@
$ws.Cells.Item(15, 54).Value = $n15_BB

$n15_BC = @'
N
'@
$ws.Cells.Item(15, 55).Value = $n15_BC

$n15_BD = $null
$ws.Cells.Item(15, 56).Value = $n15_BD

$n15_BE = $null
$ws.Cells.Item(15, 57).Value = $n15_BE

$n15_BF = $null
$ws.Cells.Item(15, 58).Value = $n15_BF

$n15_BG = @'
N
'@
$ws.Cells.Item(15, 59).Value = $n15_BG

$n15_BH = @'
A
'@
$ws.Cells.Item(15, 60).Value = $n15_BH

$n15_BI = $null
$ws.Cells.Item(15, 61).Value = $n15_BI

$n15_BJ = 0
$ws.Cells.Item(15, 62).Value = $n15_BJ

$n15_BK = $null
$ws.Cells.Item(15, 63).Value = $n15_BK

$n15_BL = @'
Bc.
'@
$ws.Cells.Item(15, 64).Value = $n15_BL

$n15_BM = @'
A
'@
$ws.Cells.Item(15, 65).Value = $n15_BM

$n15_BN = @'
None
'@
$ws.Cells.Item(15, 66).Value = $n15_BN

$n15_BO = $null
$ws.Cells.Item(15, 67).Value = $n15_BO

$n15_BP = $null
$ws.Cells.Item(15, 68).Value = $n15_BP

$n15_BQ = $null
$ws.Cells.Item(15, 69).Value = $n15_BQ

